# Apply scheduled-runner price/profit refresh to the Kujata_Profits sheets.
# Values below mirror the latest market-board pull (see commit diff for deltas).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 391.25
$ws.Range("I33").Value = 400.8
$ws.Range("K33").Value = 400.8
$ws.Range("M33").Value = -171.8
$ws.Range("H64").Value = 4117.7144
$ws.Range("I64").Value = 4366
$ws.Range("J64").Value = 3786.6667
$ws.Range("K64").Value = 4366
$ws.Range("L64").Value = 3786.6667
$ws.Range("M64").Value = -4118
$ws.Range("N64").Value = -4282.6667
$ws.Range("H67").Value = 4117.7144
$ws.Range("I67").Value = 4366
$ws.Range("J67").Value = 3786.6667
$ws.Range("K67").Value = 4366
$ws.Range("L67").Value = 3786.6667
$ws.Range("M67").Value = -3508
$ws.Range("N67").Value = -5502.6667
$ws.Range("H74").Value = 8756
$ws.Range("I74").Value = 9556.571
$ws.Range("J74").Value = 5954
$ws.Range("K74").Value = 9556.571
$ws.Range("L74").Value = 5954
$ws.Range("M74").Value = -8620.571
$ws.Range("N74").Value = -7826
$ws.Range("H77").Value = 8756
$ws.Range("I77").Value = 9556.571
$ws.Range("J77").Value = 5954
$ws.Range("K77").Value = 47782.855
$ws.Range("L77").Value = 29770
$ws.Range("M77").Value = -43102.855
$ws.Range("N77").Value = -39130
$ws.Range("H137").Value = 1135.963
$ws.Range("I137").Value = 1098.6086
$ws.Range("J137").Value = 1350.75
$ws.Range("K137").Value = 3295.8258
$ws.Range("L137").Value = 4052.25
$ws.Range("M137").Value = -745.8258000000001
$ws.Range("N137").Value = -9152.25
$ws.Range("H138").Value = 1264.4082
$ws.Range("I138").Value = 786.48
$ws.Range("J138").Value = 1762.25
$ws.Range("K138").Value = 2359.44
$ws.Range("L138").Value = 5286.75
$ws.Range("M138").Value = 2780.56
$ws.Range("N138").Value = -15566.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1469.6
$ws.Range("I74").Value = 962
$ws.Range("K74").Value = 962
$ws.Range("M74").Value = -88
$ws.Range("H77").Value = 1469.6
$ws.Range("I77").Value = 962
$ws.Range("K77").Value = 4810
$ws.Range("M77").Value = -442
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()  # leve no longer has an HQ profit figure
$ws.Range("H110").Value = 1988.3572
$ws.Range("I110").Value = 1361.8572
$ws.Range("J110").Value = 2614.8572
$ws.Range("K110").Value = 1361.8572
$ws.Range("L110").Value = 2614.8572
$ws.Range("M110").Value = 683.1428000000001
$ws.Range("N110").Value = -6704.8572

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 22728318
$ws.Range("I94").Value = 25000950
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 25000950
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -25000499
$ws.Range("N94").Value = -2902
$ws.Range("H134").Value = 3861.0557
$ws.Range("I134").Value = 1091.6207
$ws.Range("J134").Value = 15334.429
$ws.Range("K134").Value = 3274.8621
$ws.Range("L134").Value = 46003.287
$ws.Range("M134").Value = -739.8620999999998
$ws.Range("N134").Value = -51073.287

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858420
$ws.Range("I16").Value = 200001180
$ws.Range("J16").Value = 1495
$ws.Range("K16").Value = 200001180
$ws.Range("L16").Value = 1495
$ws.Range("M16").Value = -200000893
$ws.Range("N16").Value = -2069
$ws.Range("H58").Value = 1319.6
$ws.Range("I58").Value = 1347.75
$ws.Range("J58").Value = 1207
$ws.Range("K58").Value = 1347.75
$ws.Range("L58").Value = 1207
$ws.Range("M58").Value = -1144.75
$ws.Range("N58").Value = -1613
$ws.Range("H113").Value = 142858420
$ws.Range("I113").Value = 200001180
$ws.Range("J113").Value = 1495
$ws.Range("K113").Value = 200001180
$ws.Range("L113").Value = 1495
$ws.Range("M113").Value = -199999010
$ws.Range("N113").Value = -5835
$ws.Range("H136").Value = 1319.6
$ws.Range("I136").Value = 1347.75
$ws.Range("J136").Value = 1207
$ws.Range("K136").Value = 4043.25
$ws.Range("L136").Value = 3621
$ws.Range("M136").Value = -1493.25
$ws.Range("N136").Value = -8721

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2199849
$ws.Range("I4").Value = 299698.34
$ws.Range("K4").Value = 899095.02
$ws.Range("M4").Value = -898983.02
$ws.Range("H5").Value = 1240.5714
$ws.Range("I5").Value = 1624.8334
$ws.Range("K5").Value = 4874.5002
$ws.Range("M5").Value = -4762.5002
$ws.Range("H40").Value = 124.4
$ws.Range("J40").Value = 72
$ws.Range("L40").Value = 288
$ws.Range("N40").Value = -426
$ws.Range("H122").Value = 948.1429000000001
$ws.Range("I122").Value = 859.25
$ws.Range("K122").Value = 7733.25
$ws.Range("M122").Value = -5283.25
$ws.Range("H135").Value = 1240.5714
$ws.Range("I135").Value = 1624.8334
$ws.Range("K135").Value = 14623.5006
$ws.Range("M135").Value = -12088.5006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3333.3333
$ws.Range("I80").Value = 1740
$ws.Range("J80").Value = 5325
$ws.Range("K80").Value = 1740
$ws.Range("L80").Value = 5325
$ws.Range("M80").Value = -742
$ws.Range("N80").Value = -7321
$ws.Range("H83").Value = 3333.3333
$ws.Range("I83").Value = 1740
$ws.Range("J83").Value = 5325
$ws.Range("K83").Value = 8700
$ws.Range("L83").Value = 26625
$ws.Range("M83").Value = -3708
$ws.Range("N83").Value = -36609
$ws.Range("H113").Value = 2215.2
$ws.Range("I113").Value = 1419.1428
$ws.Range("J113").Value = 2643.8462
$ws.Range("K113").Value = 1419.1428
$ws.Range("L113").Value = 2643.8462
$ws.Range("M113").Value = 750.8571999999999
$ws.Range("N113").Value = -6983.8462
$ws.Range("H122").Value = 2761.4
$ws.Range("I122").Value = 2703.5
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 8110.5
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -5660.5
$ws.Range("N122").Value = -13300
$ws.Range("H124").Value = 30645
$ws.Range("J124").Value = 30645
$ws.Range("L124").Value = 30645
$ws.Range("N124").Value = -40465
$ws.Range("H126").Value = 2660.7273
$ws.Range("I126").Value = 2356
$ws.Range("J126").Value = 2728.4443
$ws.Range("K126").Value = 7068
$ws.Range("L126").Value = 8185.3329
$ws.Range("M126").Value = -4598
$ws.Range("N126").Value = -13125.3329
$ws.Range("H132").Value = 2363.85
$ws.Range("I132").Value = 2071.0557
$ws.Range("K132").Value = 6213.1671
$ws.Range("M132").Value = -3683.1671

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1258.9286
$ws.Range("J16").Value = 1516
$ws.Range("L16").Value = 1516
$ws.Range("N16").Value = -1856
$ws.Range("H61").Value = 1283.1818
$ws.Range("I61").Value = 1458.7142
$ws.Range("J61").Value = 976
$ws.Range("K61").Value = 1458.7142
$ws.Range("L61").Value = 976
$ws.Range("M61").Value = -1256.7142
$ws.Range("N61").Value = -1380
$ws.Range("H93").Value = 717.25
$ws.Range("I93").Value = 717.25
$ws.Range("K93").Value = 717.25
$ws.Range("M93").Value = 530.75
$ws.Range("H113").Value = 1283.1818
$ws.Range("I113").Value = 1458.7142
$ws.Range("J113").Value = 976
$ws.Range("K113").Value = 1458.7142
$ws.Range("L113").Value = 976
$ws.Range("M113").Value = 711.2858000000001
$ws.Range("N113").Value = -5316

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 811.25
$ws.Range("J113").Value = 1199
$ws.Range("L113").Value = 3597
$ws.Range("N113").Value = -7937
$ws.Range("H122").Value = 11306727
$ws.Range("I122").Value = 14447208
$ws.Range("J122").Value = 997.8
$ws.Range("K122").Value = 43341624
$ws.Range("L122").Value = 2993.4
$ws.Range("M122").Value = -43339174
$ws.Range("N122").Value = -7893.4
$ws.Range("H126").Value = 76924010
$ws.Range("I126").Value = 90909736
$ws.Range("K126").Value = 272729208
$ws.Range("M126").Value = -272726738
